$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-19.37033620770523,2.454455150432258,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,2.300225392562882,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,2.160312543399281,2.862205634647758,-19.37033620770523,2.554066400896631,-19.37033620770523,1.805803454192886,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,1.012605441914215,-19.37033620770523,-19.37033620770523,-19.37033620770523,2.360137661062417,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(2.963785650321165,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,2.894599130665334,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(3.608606893636345,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,1.725581646765238,-19.37033620770523,2.233232639991834),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,2.018911221869267,-19.37033620770523,2.564417559586429,-19.37033620770523,-19.37033620770523,-19.37033620770523,1.444572431067016),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,1.67426669390668,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,4.321926072584372,1.810836071890428),
    @(-19.37033620770523,-19.37033620770523,1.626626134685664,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,2.028668886408968),
    @(-19.37033620770523,-19.37033620770523,-0.2508453703785975,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,0.7365602254112755,-0.01038723039795694,-19.37033620770523,-19.37033620770523,-19.37033620770523,0.5408055831581171,0.6970574485272584,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,-19.37033620770523,0.7282318192550018,1.044659671171756,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,-19.37033620770523,1.774337592483979,-19.37033620770523,-19.37033620770523,-19.37033620770523,1.784494877832108,2.023997211490935,-19.37033620770523,-19.37033620770523),
    @(-19.37033620770523,1.594786455229055,2.133162015409021,-19.37033620770523,3.820393496930146,-19.37033620770523,2.125874410704588,2.005878757810714,-19.37033620770523,2.318332862804189),
    @(-19.37033620770523,1.741813730156444,-19.37033620770523,2.411224017401821,-19.37033620770523,3.161987610701499,2.480961814284932,-19.37033620770523,-19.37033620770523,-19.37033620770523)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

Write-Host "applied"